$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.245944912280777
$ws.Range("C2").Value = 0.2115104654614868
$ws.Range("E2").Value = 0.1764712512158786
$ws.Range("F2").Value = 2.144470902196545
$ws.Range("G2").Value = 1.093480589025646
$ws.Range("H2").Value = 1.078445673552736
$ws.Range("I2").Value = 1.055118632153878
$ws.Range("J2").Value = 0.06192150132802254
$ws.Range("L2").Value = 0.4369407038987703
$ws.Range("M2").Value = 0.3454865264307969

$ws.Range("B3").Value = 1.147694658956254
$ws.Range("C3").Value = 0.1926777248806957
$ws.Range("E3").Value = 0.1769490190079068
$ws.Range("F3").Value = 2.14837885819442
$ws.Range("G3").Value = 1.094245503618083
$ws.Range("H3").Value = 1.085045734096568
$ws.Range("I3").Value = 1.065753648113947
$ws.Range("J3").Value = 0.06162406455125691
$ws.Range("L3").Value = 0.4322979595700076
$ws.Range("M3").Value = 0.3289107360576722

$ws.Range("B4").Value = 1.087633681460147
$ws.Range("C4").Value = 0.1810279432754101
$ws.Range("E4").Value = 0.1772649399390898
$ws.Range("F4").Value = 2.152007338053309
$ws.Range("G4").Value = 1.095599167357662
$ws.Range("H4").Value = 1.089726594786811
$ws.Range("I4").Value = 1.072956982468977
$ws.Range("J4").Value = 0.06145274349237084
$ws.Range("L4").Value = 0.4295851261328352
$ws.Range("M4").Value = 0.3188395921069684

$ws.Range("B5").Value = 1.063226100100792
$ws.Range("C5").Value = 0.1762588323453542
$ws.Range("E5").Value = 0.1773993697526692
$ws.Range("F5").Value = 2.153794839951729
$ws.Range("G5").Value = 1.096372438982414
$ws.Range("H5").Value = 1.091791930362902
$ws.Range("I5").Value = 1.076061466409257
$ws.Range("J5").Value = 0.06138579498266594
$ws.Range("L5").Value = 0.4285144327189485
$ws.Range("M5").Value = 0.3147625227603257

$ws.Range("B6").Value = 1.059177364776815
$ws.Range("C6").Value = 0.1754656126184102
$ws.Range("E6").Value = 0.1774220357688701
$ws.Range("F6").Value = 2.154110303143369
$ws.Range("G6").Value = 1.09651420787516
$ws.Range("H6").Value = 1.0921444064019
$ws.Range("I6").Value = 1.076587167113217
$ws.Range("J6").Value = 0.06137485203392501
$ws.Range("L6").Value = 0.4283387516522126
$ws.Range("M6").Value = 0.3140871663987923

$ws.Range("B7").Value = 1.0873042364762
$ws.Range("C7").Value = 0.1809637133633544
$ws.Range("E7").Value = 0.1772667298507393
$ws.Range("F7").Value = 2.152030194569718
$ws.Range("G7").Value = 1.095608699409283
$ws.Range("H7").Value = 1.089753809737985
$ws.Range("I7").Value = 1.072998166434434
$ws.Range("J7").Value = 0.06145182896572621
$ws.Range("L7").Value = 0.4295705452488789
$ws.Range("M7").Value = 0.3187844976752672

$ws.Range("B8").Value = 1.212013831236334
$ws.Range("C8").Value = 0.2050348612247319
$ws.Range("E8").Value = 0.1766313097216337
$ws.Range("F8").Value = 2.145563183576456
$ws.Range("G8").Value = 1.093560464190588
$ws.Range("H8").Value = 1.080590847966619
$ws.Range("I8").Value = 1.058645665025729
$ws.Range("J8").Value = 0.06181661208349354
$ws.Range("L8").Value = 0.4353113662216117
$ws.Range("M8").Value = 0.3397492321350839

$ws.Range("B9").Value = 1.458637477956245
$ws.Range("C9").Value = 0.2515551763940493
$ws.Range("E9").Value = 0.1755637075962051
$ws.Range("F9").Value = 2.142642890364613
$ws.Range("G9").Value = 1.096588802629697
$ws.Range("H9").Value = 1.067616598355556
$ws.Range("I9").Value = 1.035855520175552
$ws.Range("J9").Value = 0.06262072777690619
$ws.Range("L9").Value = 0.4476570851095829
$ws.Range("M9").Value = 0.3816978725572184

$ws.Range("B10").Value = 1.641066206448613
$ws.Range("C10").Value = 0.2853238347464355
$ws.Range("E10").Value = 0.1748872879883203
$ws.Range("F10").Value = 2.146465728603076
$ws.Range("G10").Value = 1.103153328215086
$ws.Range("H10").Value = 1.06114140867345
$ws.Range("I10").Value = 1.022393023876496
$ws.Range("J10").Value = 0.0632644999742098
$ws.Range("L10").Value = 0.45738456243933
$ws.Range("M10").Value = 0.4130207625520654

$ws.Range("B11").Value = 1.724321643113115
$ws.Range("C11").Value = 0.3005987398458103
$ws.Range("E11").Value = 0.1746028304280056
$ws.Range("F11").Value = 2.149504835361469
$ws.Range("G11").Value = 1.107091798851428
$ws.Range("H11").Value = 1.05886221754011
$ws.Range("I11").Value = 1.016984774780845
$ws.Range("J11").Value = 0.06356864982774368
$ws.Range("L11").Value = 0.4619513168712786
$ws.Range("M11").Value = 0.4273784253054984

$ws.Range("B12").Value = 1.755886049553339
$ws.Range("C12").Value = 0.3063705569882416
$ws.Range("E12").Value = 0.1744984431146808
$ws.Range("F12").Value = 2.150842882386186
$ws.Range("G12").Value = 1.108720877202657
$ws.Range("H12").Value = 1.058095189191164
$ws.Range("I12").Value = 1.01504006114201
$ws.Range("J12").Value = 0.06368542776117181
$ws.Range("L12").Value = 0.4637008741530195
$ws.Range("M12").Value = 0.4328307600802361

$ws.Range("B13").Value = 1.74908644495639
$ws.Range("C13").Value = 0.3051280469550193
$ws.Range("E13").Value = 0.1745207768607946
$ws.Range("F13").Value = 2.150546380059183
$ws.Range("G13").Value = 1.108363891639513
$ws.Range("H13").Value = 1.058256106963881
$ws.Range("I13").Value = 1.015454291726691
$ws.Range("J13").Value = 0.0636602065824512
$ws.Range("L13").Value = 0.4633231784272596
$ws.Range("M13").Value = 0.4316558210385892

$ws.Range("B14").Value = 1.726917723006068
$ws.Range("C14").Value = 0.3010738403885966
$ws.Range("E14").Value = 0.1745941757425218
$ws.Range("F14").Value = 2.149611164009585
$ws.Range("G14").Value = 1.107223060886611
$ws.Range("H14").Value = 1.058797187349342
$ws.Range("I14").Value = 1.016822710405918
$ws.Range("J14").Value = 0.06357822521204426
$ws.Range("L14").Value = 0.4620948494664106
$ws.Range("M14").Value = 0.4278266848732386

$ws.Range("B15").Value = 1.713343575845897
$ws.Range("C15").Value = 0.2985888994261359
$ws.Range("E15").Value = 0.174639568046727
$ws.Range("F15").Value = 2.149062704098199
$ws.Range("G15").Value = 1.106542218967917
$ws.Range("H15").Value = 1.059141130151758
$ws.Range("I15").Value = 1.017674365232779
$ws.Range("J15").Value = 0.06352821739174885
$ws.Range("L15").Value = 0.4613450922748541
$ws.Range("M15").Value = 0.4254832269645803

$ws.Range("B16").Value = 1.635630558848277
$ws.Range("C16").Value = 0.2843238404733484
$ws.Range("E16").Value = 0.1749063447099477
$ws.Range("F16").Value = 2.146293298394014
$ws.Range("G16").Value = 1.102915164102228
$ws.Range("H16").Value = 1.061303788041073
$ws.Range("I16").Value = 1.022760899446183
$ws.Range("J16").Value = 0.06324484852163792
$ws.Range("L16").Value = 0.4570889526424082
$ws.Range("M16").Value = 0.4120846231695339

$ws.Range("B17").Value = 1.588023903917133
$ws.Range("C17").Value = 0.2755505000855862
$ws.Range("E17").Value = 0.175075949194845
$ws.Range("F17").Value = 2.144927519240895
$ws.Range("G17").Value = 1.100934496375928
$ws.Range("H17").Value = 1.062801349739942
$ws.Range("I17").Value = 1.026064916591636
$ws.Range("J17").Value = 0.06307388850490625
$ws.Range("L17").Value = 0.454514130940197
$ws.Range("M17").Value = 0.4038926922418185

$ws.Range("B18").Value = 1.56066707715479
$ws.Range("C18").Value = 0.270496152920856
$ws.Range("E18").Value = 0.1751756902714803
$ws.Range("F18").Value = 2.144264313529078
$ws.Range("G18").Value = 1.099884853602177
$ws.Range("H18").Value = 1.063725419371977
$ws.Range("I18").Value = 1.028032664961884
$ws.Range("J18").Value = 0.06297662116266523
$ws.Range("L18").Value = 0.4530465055041901
$ws.Range("M18").Value = 0.3991911575917086

$ws.Range("B19").Value = 1.551408898699435
$ws.Range("C19").Value = 0.2687834374596321
$ws.Range("E19").Value = 0.1752098372456432
$ws.Range("F19").Value = 2.144060770236237
$ws.Range("G19").Value = 1.099544826512613
$ws.Range("H19").Value = 1.064049057486855
$ws.Range("I19").Value = 1.028710470845773
$ws.Range("J19").Value = 0.06294387156719594
$ws.Range("L19").Value = 0.452551889057716
$ws.Range("M19").Value = 0.3976010659404494

$ws.Range("B20").Value = 1.593089107213473
$ws.Range("C20").Value = 0.2764852810545051
$ws.Range("E20").Value = 0.1750576680403615
$ws.Range("F20").Value = 2.145060243707945
$ws.Range("G20").Value = 1.101136064808145
$ws.Range("H20").Value = 1.062635439686787
$ws.Range("I20").Value = 1.025706223924836
$ws.Range("J20").Value = 0.06309197748984019
$ws.Range("L20").Value = 0.4547868451269608
$ws.Range("M20").Value = 0.4047636784224622

$ws.Range("B21").Value = 1.733428208345288
$ws.Range("C21").Value = 0.3022649972341753
$ws.Range("E21").Value = 0.1745725264201754
$ws.Range("F21").Value = 2.149880776967208
$ws.Range("G21").Value = 1.107554408486891
$ws.Range("H21").Value = 1.058635650316631
$ws.Range("I21").Value = 1.016417967268467
$ws.Range("J21").Value = 0.06360226178964012
$ws.Range("L21").Value = 0.4624550916764179
$ws.Range("M21").Value = 0.4289509789016535

$ws.Range("B22").Value = 1.825365295399251
$ws.Range("C22").Value = 0.3190409190157482
$ws.Range("E22").Value = 0.174274865813401
$ws.Range("F22").Value = 2.154122556026152
$ws.Range("G22").Value = 1.112551925326343
$ws.Range("H22").Value = 1.056581525098338
$ws.Range("I22").Value = 1.010949660039472
$ws.Range("J22").Value = 0.06394509602888832
$ws.Range("L22").Value = 0.4675845625497033
$ws.Range("M22").Value = 0.4448484313541101

$ws.Range("B23").Value = 1.776277215782329
$ws.Range("C23").Value = 0.3100939376503504
$ws.Range("E23").Value = 0.1744319612213978
$ws.Range("F23").Value = 2.151758702429206
$ws.Range("G23").Value = 1.109810961245458
$ws.Range("H23").Value = 1.057626539084481
$ws.Range("I23").Value = 1.013812997899976
$ws.Range("J23").Value = 0.0637612717517051
$ws.Range("L23").Value = 0.4648361341162257
$ws.Range("M23").Value = 0.4363555352269941

$ws.Range("B24").Value = 1.590799088622646
$ws.Range("C24").Value = 0.2760626990945241
$ws.Range("E24").Value = 0.1750659259940983
$ws.Range("F24").Value = 2.144999858968063
$ws.Range("G24").Value = 1.101044658323104
$ws.Range("H24").Value = 1.062710251071366
$ws.Range("I24").Value = 1.025868176378999
$ws.Range("J24").Value = 0.06308379628053373
$ws.Range("L24").Value = 0.4546635115564897
$ws.Range("M24").Value = 0.404369880274011

$ws.Range("B25").Value = 1.391700627520208
$ws.Range("C25").Value = 0.2390426555337228
$ws.Range("E25").Value = 0.1758335051149293
$ws.Range("F25").Value = 2.14238606390694
$ws.Range("G25").Value = 1.095010715996452
$ws.Range("H25").Value = 1.070590531867268
$ws.Range("I25").Value = 1.041445823662123
$ws.Range("J25").Value = 0.0623938079811559
$ws.Range("L25").Value = 0.4442014093590814
$ws.Range("M25").Value = 0.3702607929957438
